$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rescatables")

# Row 9 previously held the "XOTLANIHUA / XOTLANIHUA / JOEL" (ECOLOGIA, 4APV, 2)
# rescatable record. That student is moved down to a new row 15 (with only
# 1 failed subject now), and row 9 is replaced by a new rescatable record.

# New row 15: the record that used to live in row 9, now with G=1.
$ws.Range("A15").Value = 19330051920355
$ws.Range("B15").Value = "XOTLANIHUA"
$ws.Range("C15").Value = "XOTLANIHUA"
$ws.Range("D15").Value = "JOEL"
$ws.Range("E15").Value = "ECOLOGÍA"
$ws.Range("F15").Value = "4APV"
$ws.Range("G15").Value = 1

# Row 9 becomes a brand new record.
$ws.Range("A9").Value = 19330051920251
$ws.Range("B9").Value = "CRUZ"
$ws.Range("C9").Value = "PALMA"
$ws.Range("D9").Value = "VALERIA"
$ws.Range("E9").Value = "ECOLOGÍA"
$ws.Range("F9").Value = "4ASV"
$ws.Range("G9").Value = 2
